# Generate Report for Handback
# The 3bcd6935-3054-4dd3-8ac5-6f943fedbf82.md file has been handed back and
# is in sync with en-US. Update its status on the Overview, zh-cn and de-de
# sheets, and stamp the Latest Handback DateTime on the locale sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Overview sheet: row 3 is the 3bcd6935...md entry (zh-cn & de-de status cols)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# zh-cn sheet: row 3 is the 3bcd6935...md entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusText
$wsZhCn.Range("G3").Value = "2016-02-25 05:51:05"

# de-de sheet: row 3 is the 3bcd6935...md entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusText
$wsDeDe.Range("G3").Value = "2016-02-25 05:51:25"
